# Stage 1: update companies data
#
# The source rows in the sheet got reshuffled: for 49 of the 51 data rows
# (rows 2-52, i.e. every row except row 3 "BDL NORTH PARTNERS LTD" and row
# 12 "PE RESOURCES LTD"), the Company Name / Company Number / Category /
# SIC Codes / SIC Description / Typical Use Case (columns A, B, H, I, J, K)
# need to be replaced with values taken from a different row while columns
# C-G (Incorporation Date, Status, Source, Date Downloaded, Time
# Discovered) stay untouched for each row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Company Number (B) and SIC Codes (I) are treated as text so
# purely-numeric values (e.g. "64209") keep leading context and do not
# get auto-coerced into numbers by Excel.
$ws.Range("B2:B52").NumberFormat = "@"
$ws.Range("I2:I52").NumberFormat = "@"

$ws.Range("A2").Value = "GANDER INVESTMENTS LTD"
$ws.Range("B2").Value = "16473515"
$ws.Range("H2").Value = "Investments"
$ws.Range("I2").Value = "68100,68209"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

$ws.Range("A4").Value = "SEVEN (HOLDCO) LIMITED"
$ws.Range("B4").Value = "16473606"
$ws.Range("H4").Value = "Other"
$ws.Range("I4").Value = "64209"
$ws.Range("J4").Value = "Activities of other holding companies n.e.c."
$ws.Range("K4").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

$ws.Range("A5").Value = "TLJ INVESTMENT LTD"
$ws.Range("B5").Value = "16473151"
$ws.Range("H5").Value = "Investments"
$ws.Range("I5").Value = "41100,55100,68100"
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""

$ws.Range("A6").Value = "THE DISLEY GROUP LTD"
$ws.Range("B6").Value = "16473398"
$ws.Range("H6").Value = "Other"
$ws.Range("I6").Value = "64209"
$ws.Range("J6").Value = "Activities of other holding companies n.e.c."
$ws.Range("K6").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

$ws.Range("A7").Value = "GAUNT CAPITAL LTD"
$ws.Range("B7").Value = "16473262"
$ws.Range("H7").Value = "Capital"
$ws.Range("I7").Value = "64209"
$ws.Range("J7").Value = "Activities of other holding companies n.e.c."
$ws.Range("K7").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

$ws.Range("A8").Value = "INTERCONTINENTAL HOLDING COMPANY LIMITED"
$ws.Range("B8").Value = "16473418"
$ws.Range("H8").Value = "Other"
$ws.Range("I8").Value = "64209"
$ws.Range("J8").Value = "Activities of other holding companies n.e.c."
$ws.Range("K8").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

$ws.Range("A9").Value = "AJ INVESTMENT AND CONSULTANCY LTD"
$ws.Range("B9").Value = "16473328"
$ws.Range("H9").Value = "Investments"
$ws.Range("I9").Value = "64306,70229"
$ws.Range("J9").Value = "Activities of real estate investment trusts"
$ws.Range("K9").Value = "UK-regulated REIT companies."

$ws.Range("A10").Value = "BRIDGEWICK PARTNERS LIMITED"
$ws.Range("B10").Value = "16473142"
$ws.Range("H10").Value = "Partners"
$ws.Range("I10").Value = "64999"
$ws.Range("J10").Value = "Financial intermediation not elsewhere classified"
$ws.Range("K10").Value = "Catch-all credit-oriented SPVs for novel lending structures."

$ws.Range("A11").Value = "MARMIMI HOLDING LIMITED"
$ws.Range("B11").Value = "16473234"
$ws.Range("H11").Value = "Other"
$ws.Range("I11").Value = "64209"
$ws.Range("J11").Value = "Activities of other holding companies n.e.c."
$ws.Range("K11").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

$ws.Range("A13").Value = "WZRD CAPITAL MANAGEMENT LTD"
$ws.Range("B13").Value = "16472721"
$ws.Range("H13").Value = "Capital"
$ws.Range("I13").Value = "82990"
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""

$ws.Range("A14").Value = "VALANTIX LTD"
$ws.Range("B14").Value = "16472671"
$ws.Range("H14").Value = "Other"
$ws.Range("I14").Value = "62012,70100,77400"
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""

$ws.Range("A15").Value = "SYNERGIES CAPITAL (HULL) LIMITED"
$ws.Range("B15").Value = "16470828"
$ws.Range("H15").Value = "Capital"
$ws.Range("I15").Value = "68100,68209"
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""

$ws.Range("A16").Value = "GINGE INVESTMENTS LIMITED"
$ws.Range("B16").Value = "16471134"
$ws.Range("H16").Value = "Investments"
$ws.Range("I16").Value = "64209"
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""

$ws.Range("A17").Value = "HORIZON LEGACY INVESTMENTS LIMITED"
$ws.Range("B17").Value = "16470870"
$ws.Range("H17").Value = "Investments"
$ws.Range("I17").Value = "68100"
$ws.Range("J17").Value = ""
$ws.Range("K17").Value = ""

$ws.Range("A18").Value = "ARGENT CAPITAL LIMITED"
$ws.Range("B18").Value = "16470176"
$ws.Range("H18").Value = "Capital"
$ws.Range("I18").Value = "62020"
$ws.Range("J18").Value = ""
$ws.Range("K18").Value = ""

$ws.Range("A19").Value = "TEQNION UK LTD"
$ws.Range("B19").Value = "16468538"
$ws.Range("H19").Value = "Other"
$ws.Range("I19").Value = "64209,70229"
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = ""

$ws.Range("A20").Value = "ARKLE GROUP LIMITED"
$ws.Range("B20").Value = "16468549"
$ws.Range("H20").Value = "Other"
$ws.Range("I20").Value = "70221"
$ws.Range("J20").Value = "Financial management (of companies and enterprises)"
$ws.Range("K20").Value = "Treasury, capital-raising and internal financial services arm."

$ws.Range("A21").Value = "LUMINOS FUND (UK)"
$ws.Range("B21").Value = "16468863"
$ws.Range("H21").Value = "Fund"
$ws.Range("I21").Value = "85590"
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""

$ws.Range("A22").Value = "INTERNATIONAL FRESHLEAD GROUP LIMITED"
$ws.Range("B22").Value = "16468246"
$ws.Range("H22").Value = "Other"
$ws.Range("I22").Value = "64209"
$ws.Range("J22").Value = ""
$ws.Range("K22").Value = ""

$ws.Range("A23").Value = "CURRICULO PARTNERS LIMITED"
$ws.Range("B23").Value = "16467759"
$ws.Range("H23").Value = "Partners"
$ws.Range("I23").Value = "58110"
$ws.Range("J23").Value = ""
$ws.Range("K23").Value = ""

$ws.Range("A24").Value = "COASTAL MANAGEMENT SERVICES LLP"
$ws.Range("B24").Value = "OC456833"
$ws.Range("H24").Value = "LLP"
$ws.Range("I24").Value = ""
$ws.Range("J24").Value = ""
$ws.Range("K24").Value = ""

$ws.Range("A25").Value = "IREVOLUTION GROUP (TOPCO) LIMITED"
$ws.Range("B25").Value = "16464981"
$ws.Range("H25").Value = "Other"
$ws.Range("I25").Value = "64209"
$ws.Range("J25").Value = ""
$ws.Range("K25").Value = ""

$ws.Range("A26").Value = "OFFBRIDGE CAPITAL LTD"
$ws.Range("B26").Value = "16464818"
$ws.Range("H26").Value = "Capital"
$ws.Range("I26").Value = "47910,73110"
$ws.Range("J26").Value = ""
$ws.Range("K26").Value = ""

$ws.Range("A27").Value = "CITIUS CONSULTING LLP"
$ws.Range("B27").Value = "OC456812"
$ws.Range("H27").Value = "LLP"
$ws.Range("I27").Value = ""
$ws.Range("J27").Value = ""
$ws.Range("K27").Value = ""

$ws.Range("A28").Value = "FINCOM SOLUTIONS LIMITED"
$ws.Range("B28").Value = "16465217"
$ws.Range("H28").Value = "Other"
$ws.Range("I28").Value = "69202,70221,70229,74909"
$ws.Range("J28").Value = "Financial management (of companies and enterprises)"
$ws.Range("K28").Value = "Treasury, capital-raising and internal financial services arm."

$ws.Range("A29").Value = "B AND G VENTURES LTD"
$ws.Range("B29").Value = "16465082"
$ws.Range("H29").Value = "Ventures"
$ws.Range("I29").Value = "47910"
$ws.Range("J29").Value = ""
$ws.Range("K29").Value = ""

$ws.Range("A30").Value = "FINQUARK FINANCIAL TECHNOLOGIES LIMITED"
$ws.Range("B30").Value = "16465227"
$ws.Range("H30").Value = "Other"
$ws.Range("I30").Value = "46450,62020,64303"
$ws.Range("J30").Value = ""
$ws.Range("K30").Value = ""

$ws.Range("A31").Value = "CAMEL PARTNERS LTD"
$ws.Range("B31").Value = "16465520"
$ws.Range("H31").Value = "Partners"
$ws.Range("I31").Value = "68209"
$ws.Range("J31").Value = ""
$ws.Range("K31").Value = ""

$ws.Range("A32").Value = "AFROSCOT VENTURES LTD"
$ws.Range("B32").Value = "16462878"
$ws.Range("H32").Value = "Ventures"
$ws.Range("I32").Value = "70229"
$ws.Range("J32").Value = ""
$ws.Range("K32").Value = ""

$ws.Range("A33").Value = "INTERA AI HOLDINGS UK LIMITED"
$ws.Range("B33").Value = "16462908"
$ws.Range("H33").Value = "Other"
$ws.Range("I33").Value = "64305,66110,66300"
$ws.Range("J33").Value = ""
$ws.Range("K33").Value = ""

$ws.Range("A34").Value = "CFBUK SECURED FINANCE LTD"
$ws.Range("B34").Value = "16462904"
$ws.Range("H34").Value = "Other"
$ws.Range("I34").Value = "64999"
$ws.Range("J34").Value = ""
$ws.Range("K34").Value = ""

$ws.Range("A35").Value = "ST GEORGE CAPITAL (LAND) LIMITED"
$ws.Range("B35").Value = "16462880"
$ws.Range("H35").Value = "Capital"
$ws.Range("I35").Value = "68100,68209"
$ws.Range("J35").Value = ""
$ws.Range("K35").Value = ""

$ws.Range("A36").Value = "DAVIDSON CAPITAL HOLDINGS LTD"
$ws.Range("B36").Value = "SC849117"
$ws.Range("H36").Value = "Capital"
$ws.Range("I36").Value = "64209"
$ws.Range("J36").Value = ""
$ws.Range("K36").Value = ""

$ws.Range("A37").Value = "DDPT HOLDINGS LIMITED"
$ws.Range("B37").Value = "16461906"
$ws.Range("H37").Value = "Other"
$ws.Range("I37").Value = "64203,64205,70221"
$ws.Range("J37").Value = "Activities of financial services holding companies; Financial management (of companies and enterprises)"
$ws.Range("K37").Value = "Holding-company SPV for portfolio-company equity stakes, co-investment vehicles, master/feeder hubs.; Treasury, capital-raising and internal financial services arm."

$ws.Range("A38").Value = "4D CAPITAL PROPCO (44) LIMITED"
$ws.Range("B38").Value = "16461269"
$ws.Range("H38").Value = "Capital"
$ws.Range("I38").Value = "64209"
$ws.Range("J38").Value = ""
$ws.Range("K38").Value = ""

$ws.Range("A39").Value = "GROUND MORTGAGE SERVICES LIMITED"
$ws.Range("B39").Value = "16461278"
$ws.Range("H39").Value = "Other"
$ws.Range("I39").Value = "64999"
$ws.Range("J39").Value = ""
$ws.Range("K39").Value = ""

$ws.Range("A40").Value = "AMPERSAND MANAGEMENT UK LTD"
$ws.Range("B40").Value = "16460662"
$ws.Range("H40").Value = "Other"
$ws.Range("I40").Value = "66300"
$ws.Range("J40").Value = ""
$ws.Range("K40").Value = ""

$ws.Range("A41").Value = "SAMVIV PARTNERS LTD"
$ws.Range("B41").Value = "16460672"
$ws.Range("H41").Value = "Partners"
$ws.Range("I41").Value = "81229"
$ws.Range("J41").Value = ""
$ws.Range("K41").Value = ""

$ws.Range("A42").Value = "THE REEL MED LLP"
$ws.Range("B42").Value = "OC456780"
$ws.Range("H42").Value = "LLP"
$ws.Range("I42").Value = ""
$ws.Range("J42").Value = ""
$ws.Range("K42").Value = ""

$ws.Range("A43").Value = "KNOTT INVESTMENTS LIMITED"
$ws.Range("B43").Value = "16458684"
$ws.Range("H43").Value = "Investments"
$ws.Range("I43").Value = "64304"
$ws.Range("J43").Value = ""
$ws.Range("K43").Value = ""

$ws.Range("A44").Value = "TAL HOLDINGS LIMITED"
$ws.Range("B44").Value = "16458565"
$ws.Range("H44").Value = "Other"
$ws.Range("I44").Value = "64209"
$ws.Range("J44").Value = ""
$ws.Range("K44").Value = ""

$ws.Range("A45").Value = "KERSLAKE HOLDINGS LIMITED"
$ws.Range("B45").Value = "16458070"
$ws.Range("H45").Value = "Other"
$ws.Range("I45").Value = "64209"
$ws.Range("J45").Value = ""
$ws.Range("K45").Value = ""

$ws.Range("A46").Value = "KC INVESTMENTS & TRADING LIMITED"
$ws.Range("B46").Value = "16456642"
$ws.Range("H46").Value = "Investments"
$ws.Range("I46").Value = "46120,46720,64304,70229"
$ws.Range("J46").Value = ""
$ws.Range("K46").Value = ""

$ws.Range("A47").Value = "JJOHN INVESTMENTS LIMITED"
$ws.Range("B47").Value = "16456276"
$ws.Range("H47").Value = "Investments"
$ws.Range("I47").Value = "68100"
$ws.Range("J47").Value = ""
$ws.Range("K47").Value = ""

$ws.Range("A48").Value = "INTERNATIONAL CREDIT SCORE LTD"
$ws.Range("B48").Value = "16456689"
$ws.Range("H48").Value = "Other"
$ws.Range("I48").Value = "64999,66220,82912"
$ws.Range("J48").Value = ""
$ws.Range("K48").Value = ""

$ws.Range("A49").Value = "ECHO VENTURES GROUP LIMITED"
$ws.Range("B49").Value = "16455744"
$ws.Range("H49").Value = "Ventures"
$ws.Range("I49").Value = "55100,68209,68320"
$ws.Range("J49").Value = ""
$ws.Range("K49").Value = ""

$ws.Range("A50").Value = "EASEDALE HOLDINGS UK LIMITED"
$ws.Range("B50").Value = "16455608"
$ws.Range("H50").Value = "Other"
$ws.Range("I50").Value = "64991"
$ws.Range("J50").Value = ""
$ws.Range("K50").Value = ""

$ws.Range("A51").Value = "ESLB INVESTMENTS LIMITED"
$ws.Range("B51").Value = "16455669"
$ws.Range("H51").Value = "Investments"
$ws.Range("I51").Value = "68100"
$ws.Range("J51").Value = ""
$ws.Range("K51").Value = ""

$ws.Range("A52").Value = "EARLY FINANCE LIMITED"
$ws.Range("B52").Value = "16455521"
$ws.Range("H52").Value = "Other"
$ws.Range("I52").Value = "69201,70221"
$ws.Range("J52").Value = "Financial management (of companies and enterprises)"
$ws.Range("K52").Value = "Treasury, capital-raising and internal financial services arm."
